# Add new client / catalog / site-code rows to code_book.xlsx
# (client_list=sheet1, catalog=sheet2, measurement_names=sheet3, site_codes=sheet4)

$wb = $excel.ActiveWorkbook

$wsClients  = $wb.Worksheets.Item("client_list")
$wsCatalog  = $wb.Worksheets.Item("catalog")
$wsSites    = $wb.Worksheets.Item("site_codes")

# --- 1. client_list: two new golf-course clients -----------------------
[void]$wsClients.Activate()
$wsClients.Range("A37").Value = 92138
$wsClients.Range("B37").Value = "Adobe Creek National Golf Course"
$wsClients.Range("A38").Value = 92128
$wsClients.Range("B38").Value = "Chipeta Golf Course"

# --- 2. catalog: new test S022 ------------------------------------------
[void]$wsCatalog.Activate()
$wsCatalog.Range("A615").Value = "S022"
$wsCatalog.Range("B615").Value = "Organic Matter % at 440C and Sand Fractions"
[void]$wsCatalog.Range("C621").Select()

# --- 3. site_codes: codes for the two golf courses above -----------------
[void]$wsSites.Activate()
$wsSites.Range("A43").Value = "Adobe Creek National Golf Course"
$wsSites.Range("B43").Value = "ADOBE"
$wsSites.Range("C43").Value = "Golf"
$wsSites.Range("A44").Value = "Chipeta Golf Course"
$wsSites.Range("B44").Value = "CHIPETA"
$wsSites.Range("C44").Value = "Golf"

# --- 4. client_list: two landscaping clients -----------------------------
[void]$wsClients.Activate()
$wsClients.Range("A39").Value = 73293
$wsClients.Range("B39").Value = "BD Landscape Contractors"
$wsClients.Range("A40").Value = 91091
$wsClients.Range("B40").Value = "Landscape Art Inc"

# --- 5. site_codes: codes for the landscaping clients ---------------------
[void]$wsSites.Activate()
$wsSites.Range("A45").Value = "BD Landscape Contractors"
$wsSites.Range("C45").Value = "Lanscape"
$wsSites.Range("B45").Value = "BDLAND"
$wsSites.Range("A46").Value = "Landscape Art Inc"
$wsSites.Range("B46").Value = "LANDART"
$wsSites.Range("C46").Value = "Landscape"

# --- 6. client_list: Town of Basalt ---------------------------------------
[void]$wsClients.Activate()
$wsClients.Range("A41").Value = 46061
$wsClients.Range("B41").Value = "Town of Basalt"
[void]$wsClients.Range("B41").Select()

# --- 7. site_codes: code for Town of Basalt -------------------------------
[void]$wsSites.Activate()
$wsSites.Range("A47").Value = "Town of Basalt"
$wsSites.Range("B47").Value = "BASALT"
$wsSites.Range("C47").Value = "Lanscape"
[void]$wsSites.Range("A56").Select()

# --- Finish back on client_list, matching the saved view state -----------
[void]$wsClients.Activate()
[void]$wsClients.Range("B41").Select()
